# Auto-generated edit script: updates market-price related columns (H-N)
# on rows across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per the scheduled market-data refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 154.03847
$ws.Cells.Item(33, 9).Value = 161.875
$ws.Cells.Item(33, 11).Value = 161.875
$ws.Cells.Item(33, 13).Value = 67.125
$ws.Cells.Item(40, 8).Value = 1000.5769
$ws.Cells.Item(40, 9).Value = 770.4211
$ws.Cells.Item(40, 11).Value = 770.4211
$ws.Cells.Item(40, 13).Value = -595.4211
$ws.Cells.Item(64, 8).Value = 2984
$ws.Cells.Item(64, 9).Value = 2817.6667
$ws.Cells.Item(64, 10).Value = 3316.6667
$ws.Cells.Item(64, 11).Value = 2817.6667
$ws.Cells.Item(64, 12).Value = 3316.6667
$ws.Cells.Item(64, 13).Value = -2569.6667
$ws.Cells.Item(64, 14).Value = -3812.6667
$ws.Cells.Item(67, 8).Value = 2984
$ws.Cells.Item(67, 9).Value = 2817.6667
$ws.Cells.Item(67, 10).Value = 3316.6667
$ws.Cells.Item(67, 11).Value = 2817.6667
$ws.Cells.Item(67, 12).Value = 3316.6667
$ws.Cells.Item(67, 13).Value = -1959.6667
$ws.Cells.Item(67, 14).Value = -5032.6667
$ws.Cells.Item(116, 8).Value = 41674000
$ws.Cells.Item(116, 9).Value = 250000000
$ws.Cells.Item(116, 10).Value = 8801.200000000001
$ws.Cells.Item(116, 11).Value = 250000000
$ws.Cells.Item(116, 12).Value = 8801.200000000001
$ws.Cells.Item(116, 13).Value = -249996558
$ws.Cells.Item(116, 14).Value = -15685.2
$ws.Cells.Item(129, 8).Value = 769.9091
$ws.Cells.Item(129, 10).Value = 800.94116
$ws.Cells.Item(129, 12).Value = 2402.82348
$ws.Cells.Item(129, 14).Value = -12402.82348
$ws.Cells.Item(131, 8).Value = 2076.926
$ws.Cells.Item(131, 9).Value = 1313.6154
$ws.Cells.Item(131, 10).Value = 2785.7144
$ws.Cells.Item(131, 11).Value = 3940.8462
$ws.Cells.Item(131, 12).Value = 8357.143199999999
$ws.Cells.Item(131, 13).Value = 1099.1538
$ws.Cells.Item(131, 14).Value = -18437.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7948.2876
$ws.Cells.Item(32, 10).Value = 15851
$ws.Cells.Item(32, 12).Value = 15851
$ws.Cells.Item(32, 14).Value = -16425
$ws.Cells.Item(45, 8).Value = 3279.8
$ws.Cells.Item(45, 9).Value = 4123.5
$ws.Cells.Item(45, 10).Value = 2717.3333
$ws.Cells.Item(45, 11).Value = 4123.5
$ws.Cells.Item(45, 12).Value = 2717.3333
$ws.Cells.Item(45, 13).Value = -3746.5
$ws.Cells.Item(45, 14).Value = -3471.3333
$ws.Cells.Item(61, 8).Value = 2729.739
$ws.Cells.Item(61, 9).Value = 3245.2593
$ws.Cells.Item(61, 10).Value = 1997.1578
$ws.Cells.Item(61, 11).Value = 3245.2593
$ws.Cells.Item(61, 12).Value = 1997.1578
$ws.Cells.Item(61, 13).Value = -3033.2593
$ws.Cells.Item(61, 14).Value = -2421.1578
$ws.Cells.Item(74, 8).Value = 45456572
$ws.Cells.Item(74, 9).Value = 71429310
$ws.Cells.Item(74, 10).Value = 4274.875
$ws.Cells.Item(74, 11).Value = 71429310
$ws.Cells.Item(74, 12).Value = 4274.875
$ws.Cells.Item(74, 13).Value = -71428436
$ws.Cells.Item(74, 14).Value = -6022.875
$ws.Cells.Item(77, 8).Value = 45456572
$ws.Cells.Item(77, 9).Value = 71429310
$ws.Cells.Item(77, 10).Value = 4274.875
$ws.Cells.Item(77, 11).Value = 357146550
$ws.Cells.Item(77, 12).Value = 21374.375
$ws.Cells.Item(77, 13).Value = -357142182
$ws.Cells.Item(77, 14).Value = -30110.375
$ws.Cells.Item(102, 8).Value = 1439.2941
$ws.Cells.Item(102, 9).Value = 1364.5333
$ws.Cells.Item(102, 11).Value = 1364.5333
$ws.Cells.Item(102, 13).Value = 257.4666999999999
$ws.Cells.Item(122, 8).Value = 1929.0625
$ws.Cells.Item(122, 9).Value = 1853.64
$ws.Cells.Item(122, 10).Value = 2198.4285
$ws.Cells.Item(122, 11).Value = 5560.92
$ws.Cells.Item(122, 12).Value = 6595.2855
$ws.Cells.Item(122, 13).Value = -3110.92
$ws.Cells.Item(122, 14).Value = -11495.2855
$ws.Cells.Item(132, 8).Value = 10773.709
$ws.Cells.Item(132, 9).Value = 1504.9535
$ws.Cells.Item(132, 11).Value = 4514.860500000001
$ws.Cells.Item(132, 13).Value = -1984.860500000001
$ws.Cells.Item(136, 8).Value = 2729.739
$ws.Cells.Item(136, 9).Value = 3245.2593
$ws.Cells.Item(136, 10).Value = 1997.1578
$ws.Cells.Item(136, 11).Value = 9735.777900000001
$ws.Cells.Item(136, 12).Value = 5991.4734
$ws.Cells.Item(136, 13).Value = -7185.777900000001
$ws.Cells.Item(136, 14).Value = -11091.4734

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(49, 8).Value = 16000
$ws.Cells.Item(49, 10).Value = 16000
$ws.Cells.Item(49, 12).Value = 16000
$ws.Cells.Item(49, 14).Value = -16478
$ws.Cells.Item(86, 8).Value = 1978.4783
$ws.Cells.Item(86, 9).Value = 1777.3077
$ws.Cells.Item(86, 10).Value = 2240
$ws.Cells.Item(86, 11).Value = 1777.3077
$ws.Cells.Item(86, 12).Value = 2240
$ws.Cells.Item(86, 13).Value = -654.3077000000001
$ws.Cells.Item(86, 14).Value = -4486
$ws.Cells.Item(89, 8).Value = 1978.4783
$ws.Cells.Item(89, 9).Value = 1777.3077
$ws.Cells.Item(89, 10).Value = 2240
$ws.Cells.Item(89, 11).Value = 8886.538500000001
$ws.Cells.Item(89, 12).Value = 11200
$ws.Cells.Item(89, 13).Value = -3270.538500000001
$ws.Cells.Item(89, 14).Value = -22432
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4363.711
$ws.Cells.Item(31, 9).Value = 2072.8572
$ws.Cells.Item(31, 10).Value = 6368.2085
$ws.Cells.Item(31, 11).Value = 2072.8572
$ws.Cells.Item(31, 12).Value = 6368.2085
$ws.Cells.Item(31, 13).Value = -1777.8572
$ws.Cells.Item(31, 14).Value = -6958.2085
$ws.Cells.Item(34, 8).Value = 4363.711
$ws.Cells.Item(34, 9).Value = 2072.8572
$ws.Cells.Item(34, 10).Value = 6368.2085
$ws.Cells.Item(34, 11).Value = 2072.8572
$ws.Cells.Item(34, 12).Value = 6368.2085
$ws.Cells.Item(34, 13).Value = -1870.8572
$ws.Cells.Item(34, 14).Value = -6772.2085
$ws.Cells.Item(58, 8).Value = 33287.688
$ws.Cells.Item(58, 9).Value = 2301.5
$ws.Cells.Item(58, 10).Value = 64273.875
$ws.Cells.Item(58, 11).Value = 2301.5
$ws.Cells.Item(58, 12).Value = 64273.875
$ws.Cells.Item(58, 13).Value = -2098.5
$ws.Cells.Item(58, 14).Value = -64679.875
$ws.Cells.Item(86, 8).Value = 15167012
$ws.Cells.Item(86, 9).Value = 2133.3333
$ws.Cells.Item(86, 11).Value = 2133.3333
$ws.Cells.Item(86, 13).Value = -1010.3333
$ws.Cells.Item(89, 8).Value = 15167012
$ws.Cells.Item(89, 9).Value = 2133.3333
$ws.Cells.Item(89, 11).Value = 10666.6665
$ws.Cells.Item(89, 13).Value = -5050.666499999999
$ws.Cells.Item(94, 8).Value = 3411.2354
$ws.Cells.Item(94, 9).Value = 2333.6667
$ws.Cells.Item(94, 10).Value = 4623.5
$ws.Cells.Item(94, 11).Value = 2333.6667
$ws.Cells.Item(94, 12).Value = 4623.5
$ws.Cells.Item(94, 13).Value = -1882.6667
$ws.Cells.Item(94, 14).Value = -5525.5
$ws.Cells.Item(136, 8).Value = 33287.688
$ws.Cells.Item(136, 9).Value = 2301.5
$ws.Cells.Item(136, 10).Value = 64273.875
$ws.Cells.Item(136, 11).Value = 6904.5
$ws.Cells.Item(136, 12).Value = 192821.625
$ws.Cells.Item(136, 13).Value = -4354.5
$ws.Cells.Item(136, 14).Value = -197921.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 124
$ws.Cells.Item(4, 9).Value = 113.6
$ws.Cells.Item(4, 10).Value = 150
$ws.Cells.Item(4, 11).Value = 340.8
$ws.Cells.Item(4, 12).Value = 450
$ws.Cells.Item(4, 13).Value = -228.8
$ws.Cells.Item(4, 14).Value = -674
$ws.Cells.Item(114, 8).Value = 2320
$ws.Cells.Item(114, 9).Value = 3200
$ws.Cells.Item(114, 10).Value = 1000
$ws.Cells.Item(114, 11).Value = 9600
$ws.Cells.Item(114, 12).Value = 3000
$ws.Cells.Item(114, 13).Value = -6346
$ws.Cells.Item(114, 14).Value = -9508
$ws.Cells.Item(131, 8).Value = 716.12
$ws.Cells.Item(131, 9).Value = 267.27274
$ws.Cells.Item(131, 10).Value = 771.5955
$ws.Cells.Item(131, 11).Value = 801.81822
$ws.Cells.Item(131, 12).Value = 2314.7865
$ws.Cells.Item(131, 13).Value = 4238.18178
$ws.Cells.Item(131, 14).Value = -12394.7865
$ws.Cells.Item(132, 8).Value = 833.625
$ws.Cells.Item(132, 9).Value = 778.1667
$ws.Cells.Item(132, 11).Value = 7003.5003
$ws.Cells.Item(132, 13).Value = -4473.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 7000
$ws.Cells.Item(5, 10).Value = 15000
$ws.Cells.Item(5, 12).Value = 15000
$ws.Cells.Item(5, 14).Value = -15224
$ws.Cells.Item(33, 8).Value = 4745
$ws.Cells.Item(33, 10).Value = 4745
$ws.Cells.Item(33, 12).Value = 4745
$ws.Cells.Item(33, 14).Value = -5249
$ws.Cells.Item(122, 8).Value = 4055.5908
$ws.Cells.Item(122, 9).Value = 3423.0715
$ws.Cells.Item(122, 10).Value = 5162.5
$ws.Cells.Item(122, 11).Value = 10269.2145
$ws.Cells.Item(122, 12).Value = 15487.5
$ws.Cells.Item(122, 13).Value = -7819.2145
$ws.Cells.Item(122, 14).Value = -20387.5
$ws.Cells.Item(132, 8).Value = 20686.633
$ws.Cells.Item(132, 9).Value = 4292.0835
$ws.Cells.Item(132, 10).Value = 86264.836
$ws.Cells.Item(132, 11).Value = 12876.2505
$ws.Cells.Item(132, 12).Value = 258794.508
$ws.Cells.Item(132, 13).Value = -10346.2505
$ws.Cells.Item(132, 14).Value = -263854.508
$ws.Cells.Item(139, 8).Value = 36751
$ws.Cells.Item(139, 10).Value = 36751
$ws.Cells.Item(139, 12).Value = 36751
$ws.Cells.Item(139, 14).Value = -47031

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 498571.44
$ws.Cells.Item(16, 8).Value = 666.1053000000001
$ws.Cells.Item(16, 9).Value = 696.38464
$ws.Cells.Item(16, 10).Value = 600.5
$ws.Cells.Item(16, 11).Value = 696.38464
$ws.Cells.Item(16, 12).Value = 600.5
$ws.Cells.Item(16, 13).Value = -526.38464
$ws.Cells.Item(16, 14).Value = -940.5
$ws.Cells.Item(122, 8).Value = 1093062.5
$ws.Cells.Item(122, 9).Value = 1637843.9
$ws.Cells.Item(122, 11).Value = 4913531.699999999
$ws.Cells.Item(122, 13).Value = -4911081.699999999
$ws.Cells.Item(124, 8).Value = 15000
$ws.Cells.Item(124, 10).Value = 15000
$ws.Cells.Item(124, 12).Value = 15000
$ws.Cells.Item(124, 14).Value = -24820

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 200022540
$ws.Cells.Item(2, 9).Value = 500002340
$ws.Cells.Item(2, 10).Value = 36001
$ws.Cells.Item(2, 11).Value = 500002340
$ws.Cells.Item(2, 12).Value = 36001
$ws.Cells.Item(2, 13).Value = -500002228
$ws.Cells.Item(2, 14).Value = -36225
$ws.Cells.Item(40, 8).Value = 9800
$ws.Cells.Item(40, 10).Value = 9800
$ws.Cells.Item(40, 12).Value = 9800
$ws.Cells.Item(40, 14).Value = -10098
$ws.Cells.Item(107, 8).Value = 95041736
$ws.Cells.Item(107, 9).Value = 166666960
$ws.Cells.Item(107, 10).Value = 9091467
$ws.Cells.Item(107, 11).Value = 500000880
$ws.Cells.Item(107, 12).Value = 27274401
$ws.Cells.Item(107, 13).Value = -499998960
$ws.Cells.Item(107, 14).Value = -27278241
$ws.Cells.Item(132, 8).Value = 1667.7693
$ws.Cells.Item(132, 9).Value = 1186.7778
$ws.Cells.Item(132, 10).Value = 2750
$ws.Cells.Item(132, 11).Value = 3560.3334
$ws.Cells.Item(132, 12).Value = 8250
$ws.Cells.Item(132, 13).Value = -1030.3334
$ws.Cells.Item(132, 14).Value = -13310
$ws.Cells.Item(136, 8).Value = 25179478
$ws.Cells.Item(136, 9).Value = 27899934
$ws.Cells.Item(136, 10).Value = 15250
$ws.Cells.Item(136, 11).Value = 83699802
$ws.Cells.Item(136, 12).Value = 45750
$ws.Cells.Item(136, 13).Value = -83697252
$ws.Cells.Item(136, 14).Value = -50850

